$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status + timestamp text updates -------------------------------------

# Overview sheet: status columns for zh-cn / de-de, plus the "latest HO
# xliff generate date" column.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 20:52:44"

# zh-cn detail sheet: Status + Latest Handoff Datetime
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 20:52:39"

# de-de detail sheet: Status + Latest Handoff Datetime
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 20:52:44"

# --- Column width updates --------------------------------------------------
# The report regenerated the "Status" columns a bit wider to fit the new
# "Ready for handoff" text (was "In Translation").
$wsOverview.Range("E1").ColumnWidth = 16.3333333333333
$wsOverview.Range("F1").ColumnWidth = 16.3333333333333
$wsZhCn.Range("C1").ColumnWidth = 16.3333333333333
$wsDeDe.Range("C1").ColumnWidth = 16.3333333333333
